$d = $word.ActiveDocument
$h = $d.Sections(1).Headers(1)
$xml = $h.Range.WordOpenXML
$bodyStart = $xml.IndexOf("<w:body>")
$bodyEnd = $xml.IndexOf("</w:body>")
$bodyInner = $xml.Substring($bodyStart+8, $bodyEnd-$bodyStart-8)
# first paragraph only: find index of "<w:p " after start, then find matching "</w:p>" before the synthetic one
$p1start = $bodyInner.IndexOf("<w:p ")
$p1end = $bodyInner.IndexOf("</w:p>") + 6
$firstPara = $bodyInner.Substring($p1start, $p1end - $p1start)
Write-Output ("firstPara length=" + $firstPara.Length)

$newBody = "<w:body>" + $firstPara + "</w:body>"
$newDoc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15">' + $newBody + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$h.Range.InsertXML($newDoc)
Write-Output "done"
